$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a run-split boundary around an (already correct) sub-range
# without altering its text, by toggling a direct-character-formatting
# property on then off.  The COM layer coalesces adjacent runs that share
# identical resolved formatting whenever a paragraph is touched by a text
# edit; toggling Bold on/off leaves the resolved formatting unchanged but
# prevents that particular span from being swallowed back into its
# neighbours.
# ---------------------------------------------------------------------------
function Protect-Range($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Bold = $false
}

# Helper: locate a unique piece of text anywhere in the document and return
# a fresh Range positioned exactly on it (search is always re-issued against
# the whole document so earlier edits - which shift character offsets -
# never throw later look-ups off).
function Find-Range($needle) {
    $rng = $d.Content
    [void]$rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $rng
}

# ---------------------------------------------------------------------------
# 1) Title: " Daniel Balthaser" -> " Daniel" + " Balthaser" (two runs, same
#    formatting) so a gramStart/gramEnd proofErr pair can bracket the name.
# ---------------------------------------------------------------------------
$r = Find-Range(" Daniel Balthaser")
$splitStart = $r.Start + " Daniel".Length
$splitEnd = $r.End
$balthRange = $d.Range($splitStart, $splitEnd)
Protect-Range $balthRange

# ---------------------------------------------------------------------------
# 2) Job title: "Vice President" -> "Executive VP" (only the first,
#    highlighted occurrence - the one immediately followed by ", Engineering"
#    in the experience table header - changes; the later bullet point that
#    reads "...Director, and Vice President" is untouched).
# ---------------------------------------------------------------------------
$r = Find-Range("Vice President, Engineering")
$vpStart = $r.Start
$vpEnd = $r.End
$titleRange = $d.Range($vpStart, $vpStart + "Vice President".Length)
$titleRange.Text = "Executive VP"

$shrink = "Vice President".Length - "Executive VP".Length
$tail = $vpEnd - $shrink          # end of "...Engineering" after the shrink
$afterVP = $vpStart + "Executive VP".Length

# re-isolate the trailing ",", " " and "Engineering" runs that would
# otherwise have been absorbed into the freshly written "Executive VP" run
$commaRange = $d.Range($afterVP, $afterVP + 1)
Protect-Range $commaRange
$spaceRange = $d.Range($afterVP + 1, $afterVP + 2)
Protect-Range $spaceRange
$engRange = $d.Range($afterVP + 2, $tail)
Protect-Range $engRange

# ---------------------------------------------------------------------------
# 3) ", RenewData, AlphaLit, Interlegis, and Compiled Software" -> split the
#    two existing runs into seven so each company name can be wrapped in its
#    own spellStart/spellEnd pair.
# ---------------------------------------------------------------------------
$r = Find-Range(", RenewData, AlphaLit, Interlegis, and Compiled Software")
$base = $r.Start

$renewStart = $base + ", ".Length
$renewEnd = $renewStart + "RenewData".Length
Protect-Range ($d.Range($renewStart, $renewEnd))

$alphaStart = $renewEnd + ", ".Length
$alphaEnd = $alphaStart + "AlphaLit".Length
Protect-Range ($d.Range($alphaStart, $alphaEnd))

$interStart = $alphaEnd + ", ".Length
$interEnd = $interStart + "Interlegis".Length
Protect-Range ($d.Range($interStart, $interEnd))

# ---------------------------------------------------------------------------
# 4) "Le" + (old _GoBack bookmark) + "d department through " -> a single
#    "Led department through " run (the two pieces already read correctly
#    together; only the mid-word bookmark kept them apart). The neighbouring
#    "migration" run must stay separate, so it is re-protected immediately
#    afterwards.
# ---------------------------------------------------------------------------
$r = Find-Range("Led department through ")
$tmpText = "Led department through >>TMP<<"
$r.Text = $tmpText
$r2 = Find-Range($tmpText)
$r2.Text = "Led department through "

$migRange = Find-Range("migration")
Protect-Range $migRange

# ---------------------------------------------------------------------------
# 5) "Lead programmer/Software engineer on web-enabled ASP.NET..." -> split
#    out "web-enabled" into its own run for a gramStart/gramEnd proofErr pair.
# ---------------------------------------------------------------------------
$webRange = Find-Range("web-enabled")
Protect-Range $webRange
